$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string bookkeeping -------------------------------------------------
# The original sheet has "extr1".."extr8" stored in B8:B15. The target layout
# inserts two new rows ("line7","line8") right after "line6" (pushing the
# "extr*" rows down to B10:B17), and the workbook's shared-string table needs
# "line7"/"line8" to land *before* "extr1".."extr8" (matching how they'd be
# added if the two new line rows were inserted ahead of the extr rows).
# Clearing the extr* cells first (so those strings become unreferenced and
# drop out of the table), writing the two new line7/line8 strings, and only
# then re-writing extr1..extr8 reproduces that ordering.
$ws.Range("B8:B15").ClearContents()

# --- New row 8: line7 -----------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- New row 9: line8 -----------------------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Shifted row 10 (was row 8): extr1 ------------------------------------
$ws.Range("A7").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Shifted row 11 (was row 9): extr2 ------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# --- Shifted row 12 (was row 10): extr3 -----------------------------------
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# --- Shifted row 13 (was row 11): extr4 -----------------------------------
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# --- Shifted row 14 (was row 12): extr5 -----------------------------------
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# --- Shifted row 15 (was row 13): extr6 -----------------------------------
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- New row 16 (was row 14): extr7 ---------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# --- New row 17 (was row 15): extr8 ---------------------------------------
$ws.Range("A14").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

Write-Output "edit applied"
